$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 1
    3 = 3
    4 = 4
    5 = 4
    6 = 2
    7 = 3
    8 = 1
    9 = 3
    10 = 3
    11 = 5
    12 = 2
    13 = 1
    14 = 4
    15 = 8
    16 = 4
    17 = 4
    18 = 6
    19 = 2
    20 = 2
    21 = 6
    22 = 6
    24 = 5
    25 = 6
    26 = 8
    27 = 5
    28 = 3
    29 = 8
    30 = 6
    32 = 7
    33 = 3
    34 = 8
    36 = 8
    37 = 8
    38 = 5
    39 = 3
    42 = 3
    43 = 6
    45 = 2
    46 = 6
    47 = 1
    48 = 4
    49 = 4
    50 = 1
    51 = 1
    52 = 3
    53 = 1
    54 = 7
    56 = 3
    57 = 7
    58 = 8
    59 = 1
    60 = 6
    61 = 7
    63 = 7
    65 = 3
    66 = 4
    67 = 3
    68 = 5
    70 = 7
    71 = 1
    72 = 8
    73 = 1
    74 = 5
    75 = 4
    77 = 7
    78 = 6
    80 = 5
    81 = 5
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 10).Value = $values[$row]
}
